# Add US and Korea adjusted data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row 20: Republic of Korea (add this string first so it lands before
# the "United States of America" string in the shared-strings table)
$ws.Range("A20").Value = "Republic of Korea"
$ws.Range("C20").Value = 100
$ws.Range("E20").Value = 0
$ws.Range("I20").Formula = "=E20/C20"

# Row 19: "United States" -> "United States of America"
$ws.Range("A19").Value = "United States of America"

# Update selection to A20
$ws.Range("A20").Select()
